$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the run containing character offset $pos into two runs
# (boundary strictly between offset $pos-1 and $pos) without leaving any
# leftover direct formatting on either side. We do this by dropping a
# throw-away bookmark exactly at that text position (which forces the
# engine to end/start a run there) and then deleting the bookmark again -
# the run split survives even though the bookmark itself is gone.
# ---------------------------------------------------------------------------
function Split-RunAt($pos) {
    $markName = "TmpSplit" + $pos
    $d.Bookmarks.Add($markName, $d.Range($pos, $pos))
    $d.Bookmarks($markName).Delete()
}

# ---------------------------------------------------------------------------
# Step 1: turn every "Economundi" into three runs "Econo" / "M" / "undi"
# (standardizing the capitalization of the brand name), left to right.
# Because the replacement text is exactly as long as the original, none of
# the character offsets shift, so we can resolve every hit up-front.
# ---------------------------------------------------------------------------
$hits = @()
$rng = $d.Content
$rng.Find.Execute("Economundi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($rng.Find.Found) {
    $hits += $rng.Start
    $rng.Collapse(0)
    $rng.Find.Execute("Economundi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

foreach ($start in $hits) {
    # Normalize the casing in place (Economundi -> EconoMundi); this keeps
    # the text length identical so later offsets remain valid.
    $word10 = $d.Range($start, $start + 10)
    $word10.Text = "EconoMundi"

    # Split into "Econo" | "M" | "undi".
    Split-RunAt ($start + 5)
    Split-RunAt ($start + 6)
}

# ---------------------------------------------------------------------------
# Step 2: relocate the "_GoBack" bookmark from the end of the requisitos
# paragraph onto the third occurrence, landing it between the "M" and
# "undi" runs.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$thirdStart = $hits[2]
$d.Bookmarks.Add("_GoBack", $d.Range($thirdStart + 6, $thirdStart + 6))
